# chore: adapt column header formatting to respective input file names
#
# The sheet holds a diff table whose left half describes the "old"
# (FV2310) message format and whose right half the "new" (FV2404)
# message format. The column headers used to carry a generic "_old" /
# "_new" suffix; they are renamed here to carry the concrete
# <formatversion> suffix instead ("_old" -> "_FV2310", "_new" ->
# "_FV2404"). Afterwards the header range is turned into a proper
# Excel Table (ListObject) and the top row is frozen so the headers
# stay visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells -------------------------------------------
# Row 1 holds the header labels in columns A:U. Any header ending in
# "_old" becomes "..._FV2310" and any header ending in "_new" becomes
# "..._FV2404"; the lone "diff" column (K) is left untouched.
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()
    if ($header -eq $null) { continue }

    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2310"
    } elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2404"
    }
}

# --- 2) Turn the range into an Excel Table ---------------------------------
$usedRange = $ws.UsedRange
$tableRange = $ws.Range("A1:U71")
$lo = $ws.ListObjects.Add(1, $tableRange, $false, 1, "")
$lo.Name = "Table1"

# --- 3) Freeze the header row ------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
